$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = "62.922.28"
$ws.Range("E2").Value = "  +2.05%  "

# Row 3: 'Ethereum'
$ws.Range("D3").Value = "3.469.99"
$ws.Range("E3").Value = "  +2.17%  "

# Row 4: 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5: 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.65"
$ws.Range("E5").Value = "  +0.06%  "

# Row 6: 'Solana'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.13"
$ws.Range("E6").Value = "  +3.67%  "

# Row 7: 'LidoStakedEther'
$ws.Range("D7").Value = "3.472.31"
$ws.Range("E7").Value = "  +2.22%  "

# Row 8: 'USDC'
$ws.Range("E8").Value = "  -0.04%  "

# Row 9: 'XRP'
$ws.Range("E9").Value = "  +1.30%  "

# Row 10: 'Toncoin'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.65"
$ws.Range("E10").Value = "  +0.39%  "

# Row 11: 'Dogecoin'
$ws.Range("E11").Value = "  +1.35%  "

# Row 12: 'Cardano'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.403"
$ws.Range("E12").Value = "  +4.44%  "

# Row 13: 'WrappedliquidstakedEther2.0'
$ws.Range("D13").Value = "4.063.86"
$ws.Range("E13").Value = "  +2.22%  "

# Row 14: 'Avalanche'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.71"
$ws.Range("E14").Value = "  +6.33%  "

# Row 15: 'TRON'
$ws.Range("E15").Value = "  +2.72%  "

# Row 16: 'WrappedEther'
$ws.Range("D16").Value = "3.464.15"
$ws.Range("E16").Value = "  +2.23%  "

# Row 17: 'ShibaInu'
$ws.Range("E17").Value = "  +0.31%  "

# Row 18: 'WrappedBTC'
$ws.Range("D18").Value = "62.956.35"
$ws.Range("E18").Value = "  +2.02%  "

# Row 19: 'Polkadot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.32"
$ws.Range("E19").Value = "  +2.93%  "

# Row 20: 'Chainlink'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.37"
$ws.Range("E20").Value = "  +5.15%  "

# Row 21: 'Uniswap'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.23"
$ws.Range("E21").Value = "  +1.16%  "

# Row 22: 'BitcoinCash'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.05"
$ws.Range("E22").Value = "  -0.27%  "

# Row 23: 'Litecoin' -> 'Polygon'
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.557"
$ws.Range("E23").Value = "  +1.48%  "

# Row 24: 'Polygon' -> 'Litecoin'
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.70"
$ws.Range("E24").Value = "  +0.08%  "

# Row 25: 'Dai'
$ws.Range("E25").Value = "  -0.03%  "

# Row 26: 'WrappedeETH'
$ws.Range("D26").Value = "3.608.78"

# Row 27: 'PEPE'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  +1.12%  "

# Row 28: 'Kaspa'
$ws.Range("E28").Value = "  -0.58%  "

# Row 29: 'RenderToken'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("E29").Value = "  +2.43%  "

# Row 30: 'Binance-PegBSC-USD'
$ws.Range("E30").Value = "  +0.31%  "

# Row 31: 'InternetComputer(DFINITY)'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.15"
$ws.Range("E31").Value = "  +1.99%  "

# Row 32: 'PancakeSwap'
$ws.Range("E32").Value = "  -1.18%  "

# Row 34: 'Fetch.AI'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.36"
$ws.Range("E34").Value = "  -2.27%  "

# Row 35: 'EthereumClassic'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.68"
$ws.Range("E35").Value = "  +1.33%  "

# Row 36: 'EnergySwap' -> 'NEARProtocol'
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.27"
$ws.Range("E36").Value = "  +3.37%  "

# Row 37: 'NEARProtocol' -> 'EnergySwap'
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "31.96"
$ws.Range("E37").Value = "  +17.65%  "

# Row 38: 'Aptos'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.05"

# Row 39: 'Monero'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.82"
$ws.Range("E39").Value = "  +0.85%  "

# Row 40: 'ImmutableX'
$ws.Range("E40").Value = "  +5.72%  "

# Row 41: 'RenzoRestakedETH'
$ws.Range("D41").Value = "3.507.45"
$ws.Range("E41").Value = "  +2.28%  "

# Row 42: 'Hedera'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0754"
$ws.Range("E42").Value = "  -1.17%  "

# Row 43: 'Mantle'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.800"
$ws.Range("E43").Value = "  +2.20%  "

# Row 44: 'OKB'
$ws.Range("E44").Value = "  -0.31%  "

# Row 45: 'Filecoin'
$ws.Range("E45").Value = "  +0.17%  "

# Row 46: 'Stacks'
$ws.Range("E46").Value = "  +2.49%  "

# Row 47: 'ONDO'
$ws.Range("E47").Value = "  +3.86%  "

# Row 48: 'Maker'
$ws.Range("D48").Value = "2.614.64"
$ws.Range("E48").Value = "  +5.52%  "

# Row 49: 'dogwifhat'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.29"
$ws.Range("E49").Value = "  +12.20%  "

# Row 50: 'InjectiveProtocol'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.95"
$ws.Range("E50").Value = "  +1.24%  "

# Row 51: 'Cosmos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.72"
$ws.Range("E51").Value = "  +1.01%  "
